$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.659.40"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "3.689.86"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "667.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.499"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.441"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000234"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").Value = "3.666.51"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "69.651.93"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "470.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "79.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "3.835.32"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  +5.29%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.165"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").Value = "3.682.11"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.935"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("E45").Value = "  +3.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000273"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("E51").Value = "  -0.11%  "
